$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace "All the Pistol weapon's nml need fix" (4 runs: "A" / "ll the
#    Pistol weapon" / "\u2019" / "s nml need fix") with a single run of text
#    "nml for RE45 have that three image resolution:512,1024,2048", keeping
#    the eastAsia-hinted run properties (sz 30 / szCs 30 / rFonts hint
#    eastAsia) that the 2nd/4th original runs already carried, and leaving
#    the surrounding "（" / "）" runs untouched.
# ---------------------------------------------------------------------------

$run1 = "A"
$run2 = "ll the Pistol weapon"
$run3 = [string][char]0x2019
$run4 = "s nml need fix"

$txt = $d.Content.Text
$startOldText = $txt.IndexOf($run1 + $run2 + $run3 + $run4)

$run1Start = $startOldText
$run1End   = $run1Start + $run1.Length
$run2Start = $run1End
$run2End   = $run2Start + $run2.Length
$run3Start = $run2End
$run3End   = $run3Start + $run3.Length
$run4Start = $run3End
$run4End   = $run4Start + $run4.Length

# Delete the trailing runs first (back to front) so earlier offsets stay valid.
$d.Range($run4Start, $run4End).Delete()
$d.Range($run3Start, $run3End).Delete()
$d.Range($run1Start, $run1End).Delete()

# What's left is just run2's text ("ll the Pistol weapon"); overwrite it with
# the new sentence - this keeps run2's original (eastAsia-hinted) formatting.
$newText = "nml for RE45 have that three image resolution:512,1024,2048"
$r2 = $d.Range($run1Start, $run1Start + $run2.Length)
$r2.Text = $newText

$newEnd = $run1Start + $newText.Length

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: delete it from its old spot (right after
#    "Volt and ") and re-add it right after the text we just inserted above
#    (i.e. right before the closing "）").
# ---------------------------------------------------------------------------

try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
    # no pre-existing "_GoBack" bookmark - nothing to remove
}

$bmRange = $d.Range($newEnd, $newEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
